# Add a "price" column (E) to the food sheet, populating a value for every
# existing data row, then leave the selection the way the author left it
# (E103:E132 — the blank rows just below the new data, ready for more rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("E1").Value = "price"

# One price per existing food row (rows 2-102)
$prices = @(1.64,2.84,8.2799999999999994,1.66,6.92,3.2,6.99,4.76,1.56,6.81,3.07,5.05,2.06,3.45,8.01,4.5199999999999996,6.22,1.38,2.98,5.44,3.58,1.28,6.89,3.1,9.77,6.27,2.78,3.5,2.63,5.76,3.61,5.46,8.8000000000000007,5.59,8.1999999999999993,5,1.83,8.85,7.41,9.23,6.09,8.89,8.81,2.86,6.7,9.66,3.54,8.4499999999999993,9.9700000000000006,1.8,6.26,9.94,3.25,5.82,2.99,5.44,8.58,6.4,6.18,9.92,6.32,3.13,1.24,2.21,2.16,6.94,6.12,3.27,9.6199999999999992,3.99,2.9,4.5599999999999996,8.18,9.91,4.57,5.62,1.88,7.47,7.35,2.56,4.1900000000000004,5.6,8.9700000000000006,6.22,9.94,2.56,3.69,0.46,5,3.08,9.5,7.56,2.59,2.86,7.4,5.93,2.6,7.41,7.92,6.63,2.96)

for ($i = 0; $i -lt $prices.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $prices[$i]
}

# Leave the same cursor/selection state recorded in the saved workbook
$ws.Range("E103:E132").Select()
